# Implement Add Permissions for SNS. Begin work on Phone Number class for SNS
# for opting in, listing numbers, etc.
#
# This adds two new AWS/SNS-related error codes to the "Error Codes" sheet:
#   412037 -> ValidationError - The input fails to satisfy the constraints
#             specified by an AWS service.
#   412038 -> InvalidParameter - A request parameter does not comply with
#             the associated constraints.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39: new "ValidationError" entry
$ws.Cells.Item(39, 1).Value = 412037
$ws.Cells.Item(39, 2).Value = "ValidationError - The input fails to satisfy the constraints specified by an AWS service."

# Row 40: new "InvalidParameter" entry
$ws.Cells.Item(40, 1).Value = 412038
$ws.Cells.Item(40, 2).Value = "InvalidParameter - A request parameter does not comply with the associated constraints."

# Move the active selection to just past the newly added data, like the
# author left it after entering the new rows.
$ws.Activate()
$ws.Range("A41").Select()
